# Append new rows to Sheet2 and Sheet3 (ThemisAutomation gold-data import).
# Order matches the source tool's interleaved sheet2/sheet3 writes so that
# newly-created shared-string entries land at the same indices as the target export.
$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws3 = $wb.Worksheets.Item("Sheet3")

# Each entry: sheet variable name, row number, value
$newRows = @(
    @('ws2', 430, '10-1307764'),
    @('ws3', 371, 'C-106761'),
    @('ws3', 372, 'OF-222433'),
    @('ws2', 431, '10-1307792'),
    @('ws3', 373, 'C-106767'),
    @('ws3', 374, 'C-117487'),
    @('ws2', 432, '10-1320426'),
    @('ws2', 433, '10-1320320'),
    @('ws2', 434, '10-1320418'),
    @('ws3', 375, 'C-117502'),
    @('ws2', 435, '10-1320419'),
    @('ws3', 376, 'C-117509'),
    @('ws2', 436, '10-1320421'),
    @('ws3', 377, 'C-117517'),
    @('ws2', 437, '10-1320424'),
    @('ws3', 378, 'C-117518'),
    @('ws2', 438, '10-1330620'),
    @('ws3', 379, 'C-117523'),
    @('ws2', 439, '10-1307793'),
    @('ws3', 380, 'C-106769'),
    @('ws2', 440, '10-1307794'),
    @('ws3', 381, 'C-106770'),
    @('ws2', 441, '10-1307795'),
    @('ws3', 382, 'C-106771'),
    @('ws2', 442, '10-1307796'),
    @('ws3', 383, 'C-106772'),
    @('ws2', 443, '10-1307798'),
    @('ws3', 384, 'C-106774'),
    @('ws2', 444, '10-1330637'),
    @('ws3', 385, 'C-117586'),
    @('ws2', 445, '10-1330628'),
    @('ws3', 386, 'C-117588'),
    @('ws2', 446, '10-1330639'),
    @('ws3', 387, 'C-117590'),
    @('ws2', 447, '10-1307799'),
    @('ws3', 388, 'C-106775'),
    @('ws2', 448, '10-1307800'),
    @('ws3', 389, 'C-106776'),
    @('ws2', 449, '10-1307801'),
    @('ws3', 390, 'C-106777'),
    @('ws2', 450, '10-1307802'),
    @('ws3', 391, 'C-106778'),
    @('ws2', 451, '10-1307803'),
    @('ws3', 392, 'C-106779'),
    @('ws2', 452, '10-1307804'),
    @('ws3', 393, 'C-106781'),
    @('ws2', 453, '10-1307805'),
    @('ws3', 394, 'C-106782'),
    @('ws2', 454, '10-1307905'),
    @('ws3', 395, 'OF-222470'),
    @('ws2', 455, '10-1307980'),
    @('ws2', 456, '10-1307981'),
    @('ws3', 396, 'C-106784'),
    @('ws2', 457, '10-1307985'),
    @('ws2', 458, '10-1307984'),
    @('ws3', 397, 'C-106785'),
    @('ws2', 459, '10-1307986'),
    @('ws3', 398, 'C-106786'),
    @('ws2', 460, '10-1307989'),
    @('ws3', 399, 'C-106787'),
    @('ws2', 461, '10-1307992'),
    @('ws3', 400, 'C-106788'),
    @('ws2', 462, '10-1308042'),
    @('ws3', 401, 'C-106789'),
    @('ws2', 463, '10-1308155'),
    @('ws3', 402, 'C-106818'),
    @('ws2', 464, '10-1308156'),
    @('ws3', 403, 'C-106819'),
    @('ws2', 465, '10-1308158'),
    @('ws3', 404, 'C-106820'),
    @('ws2', 466, '10-1308161'),
    @('ws3', 405, 'C-106826'),
    @('ws2', 467, '10-1308164'),
    @('ws3', 406, 'C-106828'),
    @('ws2', 468, '10-1308165'),
    @('ws3', 407, 'C-106829'),
    @('ws2', 469, '10-1308167'),
    @('ws3', 408, 'C-106831'),
    @('ws2', 470, '10-1308168'),
    @('ws3', 409, 'C-106832')
)

foreach ($entry in $newRows) {
    $targetSheet = if ($entry[0] -eq "ws2") { $ws2 } else { $ws3 }
    $targetSheet.Cells.Item($entry[1], 1).Value = $entry[2]
}

# Restore view state: Sheet3 selection, then Sheet2 (kept as the active tab, matching the workbook).
[void]$ws3.Range("A390").Select()
[void]$ws2.Range("A470").Select()
